# Adds three new rows (5,6,7) of script-line data to the "Политод" sheet
# and extends row 4 (previously only holding a filename cell) with blank
# but styled B:E cells, matching the upstream "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: reuse the bordered style from row 3 (s=8/8/9/9/9), keep the
#     existing filename in A4, leave B4:E4 blank (just styled). -----------
$ws.Range("A3:E3").Copy($ws.Range("A4:E4"))
$ws.Cells.Item(4, 1).Value = "SCRIPT/T01P02A/um0312.ssb"
$ws.Range("B4:E4").ClearContents()

# --- Rows 5 & 6: same bordered style pattern as row 3/4. ------------------
$ws.Range("A3:E3").Copy($ws.Range("A5:E5"))
$ws.Range("A3:E3").Copy($ws.Range("A6:E6"))

# Row heights match the rest of the data rows.
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(5).RowHeight = 43.2
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 43.2

# Row 5 content (filename, English, RU translation, "converted" RU text).
$ws.Cells.Item(5, 1).Value = "SCRIPT/T01P02A/um0602.ssb"
$ws.Cells.Item(5, 3).Value = ' Ha ha ha!'
$ws.Cells.Item(5, 4).Value = ' Ха-ха-ха!'
$ws.Cells.Item(5, 5).Value = ' Öà-öà-öà!'
$ws.Cells.Item(5, 2).Value = 332

# Row 6 content.
$ws.Cells.Item(6, 3).Value = ' I wonder if there\''s something\nspecial if you recycle a lot…'
$ws.Cells.Item(6, 1).Value = "SCRIPT/P01P04A/um0717.ssb"
$ws.Cells.Item(6, 4).Value = ' Интересно, если переработать\nмного вещей, произойдёт ли что-то\nособенное…'
$ws.Cells.Item(6, 5).Value = ' Éîóåñåòîï, åòìé ðåñåñàáïóàóû\níîïãï âåþåê, ðñïéèïêäæó ìé œóï-óï\nïòïáåîîïå...'
$ws.Cells.Item(6, 2).Value = 307

# --- Row 7: final row of the table, so it uses the plain (no-border)
#     column styles (s=4/4/5/5/5), same as the old trailing row used. -----
$ws.Cells.Item(7, 3).Value = ' An expedition! ♪[K] I envy you! ♪'
$ws.Cells.Item(7, 1).Value = "SCRIPT/G01P03A/um0803.ssb"
$ws.Cells.Item(7, 4).Value = ' Экспедиция! ♪[K] Я вам завидую! ♪'
$ws.Cells.Item(7, 5).Value = ' Üëòðåäéøéÿ! ♪[K] Ÿ âàí èàâéäôý! ♪'
$ws.Cells.Item(7, 2).Value = 288

# Match the updated selection from the saved view state.
$ws.Range("C4").Select()
